$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-85 (permuted data) - only columns D, J, K, L, M, N, O, P change
$ws.Cells.Item(2, 4).Value = 45096
$ws.Cells.Item(2, 10).Value = 65
$ws.Cells.Item(2, 11).Value = 31000
$ws.Cells.Item(2, 12).Value = 32000
$ws.Cells.Item(2, 13).Value = 31538
$ws.Cells.Item(2, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(2, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(2, 16).Value = 1262
$ws.Cells.Item(3, 4).Value = 44767
$ws.Cells.Item(3, 10).Value = 45
$ws.Cells.Item(3, 11).Value = 37000
$ws.Cells.Item(3, 12).Value = 38000
$ws.Cells.Item(3, 13).Value = 37556
$ws.Cells.Item(3, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(3, 16).Value = 1502
$ws.Cells.Item(4, 4).Value = 44908
$ws.Cells.Item(4, 10).Value = 20
$ws.Cells.Item(4, 11).Value = 18000
$ws.Cells.Item(4, 12).Value = 18000
$ws.Cells.Item(4, 13).Value = 18000
$ws.Cells.Item(4, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(4, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 16).Value = 720
$ws.Cells.Item(5, 4).Value = 44679
$ws.Cells.Item(5, 10).Value = 77
$ws.Cells.Item(5, 11).Value = 26000
$ws.Cells.Item(5, 12).Value = 27000
$ws.Cells.Item(5, 13).Value = 26506
$ws.Cells.Item(5, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(5, 16).Value = 1060
$ws.Cells.Item(6, 4).Value = 44510
$ws.Cells.Item(6, 10).Value = 73
$ws.Cells.Item(6, 11).Value = 16500
$ws.Cells.Item(6, 12).Value = 17000
$ws.Cells.Item(6, 13).Value = 16740
$ws.Cells.Item(6, 16).Value = 670
$ws.Cells.Item(7, 4).Value = 44543
$ws.Cells.Item(7, 10).Value = 35
$ws.Cells.Item(7, 11).Value = 18000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 18000
$ws.Cells.Item(7, 16).Value = 720
$ws.Cells.Item(8, 4).Value = 44847
$ws.Cells.Item(8, 10).Value = 71
$ws.Cells.Item(8, 11).Value = 30000
$ws.Cells.Item(8, 12).Value = 31000
$ws.Cells.Item(8, 13).Value = 30493
$ws.Cells.Item(8, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 1220
$ws.Cells.Item(9, 4).Value = 44165
$ws.Cells.Item(9, 10).Value = 45
$ws.Cells.Item(9, 11).Value = 22000
$ws.Cells.Item(9, 12).Value = 22000
$ws.Cells.Item(9, 13).Value = 22000
$ws.Cells.Item(9, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(9, 16).Value = 880
$ws.Cells.Item(10, 4).Value = 44484
$ws.Cells.Item(10, 10).Value = 71
$ws.Cells.Item(10, 12).Value = 30000
$ws.Cells.Item(10, 13).Value = 29507
$ws.Cells.Item(10, 16).Value = 1180
$ws.Cells.Item(11, 4).Value = 44250
$ws.Cells.Item(11, 10).Value = 38
$ws.Cells.Item(11, 11).Value = 18000
$ws.Cells.Item(11, 12).Value = 18000
$ws.Cells.Item(11, 13).Value = 18000
$ws.Cells.Item(11, 15).Value = "Provincia de Talca"
$ws.Cells.Item(11, 16).Value = 720
$ws.Cells.Item(12, 4).Value = 44453
$ws.Cells.Item(12, 11).Value = 21000
$ws.Cells.Item(12, 12).Value = 22000
$ws.Cells.Item(12, 13).Value = 21521
$ws.Cells.Item(12, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(12, 16).Value = 861
$ws.Cells.Item(13, 4).Value = 44858
$ws.Cells.Item(13, 10).Value = 80
$ws.Cells.Item(13, 11).Value = 24000
$ws.Cells.Item(13, 12).Value = 25000
$ws.Cells.Item(13, 13).Value = 24500
$ws.Cells.Item(13, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 980
$ws.Cells.Item(14, 4).Value = 44411
$ws.Cells.Item(14, 10).Value = 35
$ws.Cells.Item(14, 11).Value = 34000
$ws.Cells.Item(14, 12).Value = 34000
$ws.Cells.Item(14, 13).Value = 34000
$ws.Cells.Item(14, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(14, 16).Value = 1360
$ws.Cells.Item(15, 4).Value = 44900
$ws.Cells.Item(15, 10).Value = 73
$ws.Cells.Item(15, 11).Value = 21000
$ws.Cells.Item(15, 12).Value = 22000
$ws.Cells.Item(15, 13).Value = 21479
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 859
$ws.Cells.Item(16, 4).Value = 44529
$ws.Cells.Item(16, 10).Value = 73
$ws.Cells.Item(16, 11).Value = 17000
$ws.Cells.Item(16, 12).Value = 18000
$ws.Cells.Item(16, 13).Value = 17521
$ws.Cells.Item(16, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(16, 16).Value = 701
$ws.Cells.Item(17, 4).Value = 45100
$ws.Cells.Item(17, 10).Value = 35
$ws.Cells.Item(17, 11).Value = 34000
$ws.Cells.Item(17, 12).Value = 34000
$ws.Cells.Item(17, 13).Value = 34000
$ws.Cells.Item(17, 16).Value = 1360
$ws.Cells.Item(18, 4).Value = 44487
$ws.Cells.Item(18, 11).Value = 20000
$ws.Cells.Item(18, 12).Value = 21000
$ws.Cells.Item(18, 13).Value = 20521
$ws.Cells.Item(18, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(18, 16).Value = 821
$ws.Cells.Item(19, 4).Value = 44848
$ws.Cells.Item(19, 11).Value = 30000
$ws.Cells.Item(19, 12).Value = 30000
$ws.Cells.Item(19, 13).Value = 30000
$ws.Cells.Item(19, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(19, 16).Value = 1200
$ws.Cells.Item(20, 4).Value = 44343
$ws.Cells.Item(20, 10).Value = 40
$ws.Cells.Item(20, 11).Value = 28000
$ws.Cells.Item(20, 12).Value = 28000
$ws.Cells.Item(20, 13).Value = 28000
$ws.Cells.Item(20, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(20, 16).Value = 1120
$ws.Cells.Item(21, 4).Value = 45043
$ws.Cells.Item(21, 10).Value = 38
$ws.Cells.Item(21, 11).Value = 29000
$ws.Cells.Item(21, 12).Value = 29000
$ws.Cells.Item(21, 13).Value = 29000
$ws.Cells.Item(21, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(21, 16).Value = 1160
$ws.Cells.Item(22, 4).Value = 44917
$ws.Cells.Item(22, 10).Value = 38
$ws.Cells.Item(22, 11).Value = 28000
$ws.Cells.Item(22, 12).Value = 28000
$ws.Cells.Item(22, 13).Value = 28000
$ws.Cells.Item(22, 15).Value = "Región Metropolitana"
$ws.Cells.Item(22, 16).Value = 1120
$ws.Cells.Item(23, 4).Value = 44469
$ws.Cells.Item(23, 10).Value = 73
$ws.Cells.Item(23, 11).Value = 28000
$ws.Cells.Item(23, 12).Value = 29000
$ws.Cells.Item(23, 13).Value = 28521
$ws.Cells.Item(23, 16).Value = 1141
$ws.Cells.Item(24, 4).Value = 44526
$ws.Cells.Item(24, 10).Value = 73
$ws.Cells.Item(24, 11).Value = 16000
$ws.Cells.Item(24, 12).Value = 17000
$ws.Cells.Item(24, 13).Value = 16521
$ws.Cells.Item(24, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(24, 16).Value = 661
$ws.Cells.Item(25, 4).Value = 44252
$ws.Cells.Item(25, 10).Value = 40
$ws.Cells.Item(25, 11).Value = 18000
$ws.Cells.Item(25, 12).Value = 19000
$ws.Cells.Item(25, 13).Value = 18625
$ws.Cells.Item(25, 15).Value = "Provincia de Talca"
$ws.Cells.Item(25, 16).Value = 745
$ws.Cells.Item(26, 4).Value = 44399
$ws.Cells.Item(26, 10).Value = 38
$ws.Cells.Item(26, 11).Value = 33000
$ws.Cells.Item(26, 12).Value = 33000
$ws.Cells.Item(26, 13).Value = 33000
$ws.Cells.Item(26, 16).Value = 1320
$ws.Cells.Item(27, 4).Value = 44372
$ws.Cells.Item(27, 10).Value = 50
$ws.Cells.Item(27, 11).Value = 33000
$ws.Cells.Item(27, 12).Value = 34000
$ws.Cells.Item(27, 13).Value = 33500
$ws.Cells.Item(27, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(27, 16).Value = 1340
$ws.Cells.Item(28, 4).Value = 44476
$ws.Cells.Item(28, 11).Value = 23000
$ws.Cells.Item(28, 12).Value = 24000
$ws.Cells.Item(28, 13).Value = 23521
$ws.Cells.Item(28, 16).Value = 941
$ws.Cells.Item(29, 4).Value = 44567
$ws.Cells.Item(29, 10).Value = 68
$ws.Cells.Item(29, 11).Value = 24000
$ws.Cells.Item(29, 12).Value = 25000
$ws.Cells.Item(29, 13).Value = 24559
$ws.Cells.Item(29, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(29, 16).Value = 982
$ws.Cells.Item(30, 4).Value = 44515
$ws.Cells.Item(30, 10).Value = 73
$ws.Cells.Item(30, 11).Value = 16000
$ws.Cells.Item(30, 12).Value = 17000
$ws.Cells.Item(30, 13).Value = 16521
$ws.Cells.Item(30, 16).Value = 661
$ws.Cells.Item(31, 4).Value = 44876
$ws.Cells.Item(31, 11).Value = 18000
$ws.Cells.Item(31, 12).Value = 18000
$ws.Cells.Item(31, 13).Value = 18000
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 720
$ws.Cells.Item(32, 4).Value = 44525
$ws.Cells.Item(32, 10).Value = 73
$ws.Cells.Item(32, 11).Value = 16000
$ws.Cells.Item(32, 12).Value = 17000
$ws.Cells.Item(32, 13).Value = 16479
$ws.Cells.Item(32, 15).Value = "Provincia de Talca"
$ws.Cells.Item(32, 16).Value = 659
$ws.Cells.Item(33, 4).Value = 44536
$ws.Cells.Item(33, 10).Value = 81
$ws.Cells.Item(33, 11).Value = 27000
$ws.Cells.Item(33, 12).Value = 28000
$ws.Cells.Item(33, 13).Value = 27556
$ws.Cells.Item(33, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(33, 16).Value = 1102
$ws.Cells.Item(34, 4).Value = 44524
$ws.Cells.Item(34, 10).Value = 65
$ws.Cells.Item(34, 11).Value = 16000
$ws.Cells.Item(34, 12).Value = 17000
$ws.Cells.Item(34, 13).Value = 16538
$ws.Cells.Item(34, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(34, 16).Value = 662
$ws.Cells.Item(35, 4).Value = 44601
$ws.Cells.Item(35, 11).Value = 23000
$ws.Cells.Item(35, 12).Value = 24000
$ws.Cells.Item(35, 13).Value = 23600
$ws.Cells.Item(35, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(35, 15).Value = "Región Metropolitana"
$ws.Cells.Item(35, 16).Value = 944
$ws.Cells.Item(36, 4).Value = 44162
$ws.Cells.Item(36, 10).Value = 35
$ws.Cells.Item(36, 11).Value = 17000
$ws.Cells.Item(36, 12).Value = 17000
$ws.Cells.Item(36, 13).Value = 17000
$ws.Cells.Item(36, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(36, 16).Value = 680
$ws.Cells.Item(37, 4).Value = 44868
$ws.Cells.Item(37, 10).Value = 76
$ws.Cells.Item(37, 11).Value = 22000
$ws.Cells.Item(37, 12).Value = 23000
$ws.Cells.Item(37, 13).Value = 22500
$ws.Cells.Item(37, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(37, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(37, 16).Value = 900
$ws.Cells.Item(38, 4).Value = 44592
$ws.Cells.Item(38, 11).Value = 22000
$ws.Cells.Item(38, 12).Value = 22000
$ws.Cells.Item(38, 13).Value = 22000
$ws.Cells.Item(38, 16).Value = 880
$ws.Cells.Item(39, 4).Value = 44159
$ws.Cells.Item(39, 10).Value = 35
$ws.Cells.Item(39, 11).Value = 22000
$ws.Cells.Item(39, 12).Value = 22000
$ws.Cells.Item(39, 13).Value = 22000
$ws.Cells.Item(39, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(39, 16).Value = 880
$ws.Cells.Item(40, 4).Value = 44181
$ws.Cells.Item(40, 10).Value = 38
$ws.Cells.Item(40, 11).Value = 26000
$ws.Cells.Item(40, 12).Value = 26000
$ws.Cells.Item(40, 13).Value = 26000
$ws.Cells.Item(40, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(40, 15).Value = "Región Metropolitana"
$ws.Cells.Item(40, 16).Value = 1040
$ws.Cells.Item(41, 4).Value = 44452
$ws.Cells.Item(41, 10).Value = 70
$ws.Cells.Item(41, 11).Value = 31000
$ws.Cells.Item(41, 12).Value = 32000
$ws.Cells.Item(41, 13).Value = 31500
$ws.Cells.Item(41, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(41, 16).Value = 1260
$ws.Cells.Item(42, 4).Value = 44523
$ws.Cells.Item(42, 10).Value = 70
$ws.Cells.Item(42, 11).Value = 16000
$ws.Cells.Item(42, 12).Value = 16500
$ws.Cells.Item(42, 13).Value = 16250
$ws.Cells.Item(42, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(42, 15).Value = "Provincia de Talca"
$ws.Cells.Item(42, 16).Value = 650
$ws.Cells.Item(43, 4).Value = 44874
$ws.Cells.Item(43, 10).Value = 78
$ws.Cells.Item(43, 11).Value = 20000
$ws.Cells.Item(43, 12).Value = 21000
$ws.Cells.Item(43, 13).Value = 20513
$ws.Cells.Item(43, 16).Value = 821
$ws.Cells.Item(44, 4).Value = 44537
$ws.Cells.Item(44, 10).Value = 78
$ws.Cells.Item(44, 11).Value = 27000
$ws.Cells.Item(44, 12).Value = 28000
$ws.Cells.Item(44, 13).Value = 27487
$ws.Cells.Item(44, 16).Value = 1099
$ws.Cells.Item(45, 4).Value = 44629
$ws.Cells.Item(45, 10).Value = 45
$ws.Cells.Item(45, 11).Value = 24000
$ws.Cells.Item(45, 12).Value = 25000
$ws.Cells.Item(45, 13).Value = 24444
$ws.Cells.Item(45, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(45, 15).Value = "Región Metropolitana"
$ws.Cells.Item(45, 16).Value = 978
$ws.Cells.Item(46, 4).Value = 44475
$ws.Cells.Item(46, 11).Value = 25000
$ws.Cells.Item(46, 12).Value = 26000
$ws.Cells.Item(46, 13).Value = 25479
$ws.Cells.Item(46, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(46, 16).Value = 1019
$ws.Cells.Item(47, 4).Value = 44910
$ws.Cells.Item(47, 10).Value = 35
$ws.Cells.Item(47, 11).Value = 19000
$ws.Cells.Item(47, 12).Value = 19000
$ws.Cells.Item(47, 13).Value = 19000
$ws.Cells.Item(47, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(47, 16).Value = 760
$ws.Cells.Item(48, 4).Value = 44448
$ws.Cells.Item(48, 10).Value = 45
$ws.Cells.Item(48, 11).Value = 32000
$ws.Cells.Item(48, 12).Value = 32000
$ws.Cells.Item(48, 13).Value = 32000
$ws.Cells.Item(48, 16).Value = 1280
$ws.Cells.Item(49, 4).Value = 44410
$ws.Cells.Item(49, 11).Value = 34000
$ws.Cells.Item(49, 12).Value = 34000
$ws.Cells.Item(49, 13).Value = 34000
$ws.Cells.Item(49, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(49, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(49, 16).Value = 1360
$ws.Cells.Item(50, 4).Value = 44370
$ws.Cells.Item(50, 10).Value = 45
$ws.Cells.Item(50, 11).Value = 32000
$ws.Cells.Item(50, 12).Value = 32000
$ws.Cells.Item(50, 13).Value = 32000
$ws.Cells.Item(50, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(50, 16).Value = 1280
$ws.Cells.Item(51, 4).Value = 44509
$ws.Cells.Item(51, 10).Value = 80
$ws.Cells.Item(51, 11).Value = 15000
$ws.Cells.Item(51, 12).Value = 16000
$ws.Cells.Item(51, 13).Value = 15500
$ws.Cells.Item(51, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(51, 16).Value = 620
$ws.Cells.Item(52, 4).Value = 44161
$ws.Cells.Item(52, 10).Value = 35
$ws.Cells.Item(52, 11).Value = 21000
$ws.Cells.Item(52, 12).Value = 21000
$ws.Cells.Item(52, 13).Value = 21000
$ws.Cells.Item(52, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(52, 16).Value = 840
$ws.Cells.Item(53, 4).Value = 44882
$ws.Cells.Item(53, 10).Value = 65
$ws.Cells.Item(53, 11).Value = 19000
$ws.Cells.Item(53, 12).Value = 20000
$ws.Cells.Item(53, 13).Value = 19462
$ws.Cells.Item(53, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(53, 15).Value = "Región Metropolitana"
$ws.Cells.Item(53, 16).Value = 778
$ws.Cells.Item(54, 4).Value = 44160
$ws.Cells.Item(54, 10).Value = 35
$ws.Cells.Item(54, 11).Value = 21000
$ws.Cells.Item(54, 13).Value = 21000
$ws.Cells.Item(54, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(54, 16).Value = 840
$ws.Cells.Item(55, 4).Value = 44546
$ws.Cells.Item(55, 10).Value = 75
$ws.Cells.Item(55, 12).Value = 18500
$ws.Cells.Item(55, 13).Value = 18267
$ws.Cells.Item(55, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(55, 16).Value = 731
$ws.Cells.Item(56, 4).Value = 44365
$ws.Cells.Item(56, 10).Value = 70
$ws.Cells.Item(56, 12).Value = 23000
$ws.Cells.Item(56, 13).Value = 22500
$ws.Cells.Item(56, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(56, 16).Value = 900
$ws.Cells.Item(57, 4).Value = 45113
$ws.Cells.Item(57, 10).Value = 60
$ws.Cells.Item(57, 11).Value = 22000
$ws.Cells.Item(57, 12).Value = 23000
$ws.Cells.Item(57, 13).Value = 22333
$ws.Cells.Item(57, 16).Value = 893
$ws.Cells.Item(58, 4).Value = 44628
$ws.Cells.Item(58, 10).Value = 73
$ws.Cells.Item(58, 11).Value = 23000
$ws.Cells.Item(58, 12).Value = 24000
$ws.Cells.Item(58, 13).Value = 23521
$ws.Cells.Item(58, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(58, 16).Value = 941
$ws.Cells.Item(59, 4).Value = 44406
$ws.Cells.Item(59, 10).Value = 35
$ws.Cells.Item(59, 11).Value = 32000
$ws.Cells.Item(59, 12).Value = 32000
$ws.Cells.Item(59, 13).Value = 32000
$ws.Cells.Item(59, 16).Value = 1280
$ws.Cells.Item(60, 4).Value = 44508
$ws.Cells.Item(60, 10).Value = 68
$ws.Cells.Item(60, 11).Value = 16000
$ws.Cells.Item(60, 12).Value = 17000
$ws.Cells.Item(60, 13).Value = 16515
$ws.Cells.Item(60, 16).Value = 661
$ws.Cells.Item(61, 4).Value = 44895
$ws.Cells.Item(61, 11).Value = 22000
$ws.Cells.Item(61, 12).Value = 23000
$ws.Cells.Item(61, 13).Value = 22521
$ws.Cells.Item(61, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(61, 15).Value = "Región Metropolitana"
$ws.Cells.Item(61, 16).Value = 901
$ws.Cells.Item(62, 4).Value = 44532
$ws.Cells.Item(62, 10).Value = 73
$ws.Cells.Item(62, 11).Value = 28000
$ws.Cells.Item(62, 12).Value = 29000
$ws.Cells.Item(62, 13).Value = 28521
$ws.Cells.Item(62, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(62, 16).Value = 1141
$ws.Cells.Item(63, 4).Value = 44894
$ws.Cells.Item(63, 10).Value = 73
$ws.Cells.Item(63, 11).Value = 22000
$ws.Cells.Item(63, 12).Value = 23000
$ws.Cells.Item(63, 13).Value = 22521
$ws.Cells.Item(63, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(63, 15).Value = "Región Metropolitana"
$ws.Cells.Item(63, 16).Value = 901
$ws.Cells.Item(64, 4).Value = 44907
$ws.Cells.Item(64, 11).Value = 18000
$ws.Cells.Item(64, 12).Value = 19000
$ws.Cells.Item(64, 13).Value = 18521
$ws.Cells.Item(64, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(64, 16).Value = 741
$ws.Cells.Item(65, 4).Value = 44550
$ws.Cells.Item(65, 11).Value = 17000
$ws.Cells.Item(65, 12).Value = 18000
$ws.Cells.Item(65, 13).Value = 17521
$ws.Cells.Item(65, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(65, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(65, 16).Value = 701
$ws.Cells.Item(66, 4).Value = 44483
$ws.Cells.Item(66, 10).Value = 55
$ws.Cells.Item(66, 11).Value = 29000
$ws.Cells.Item(66, 12).Value = 30000
$ws.Cells.Item(66, 13).Value = 29455
$ws.Cells.Item(66, 16).Value = 1178
$ws.Cells.Item(67, 4).Value = 44875
$ws.Cells.Item(67, 10).Value = 68
$ws.Cells.Item(67, 11).Value = 15000
$ws.Cells.Item(67, 12).Value = 16000
$ws.Cells.Item(67, 13).Value = 15559
$ws.Cells.Item(67, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(67, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(67, 16).Value = 622
$ws.Cells.Item(68, 4).Value = 44831
$ws.Cells.Item(68, 10).Value = 45
$ws.Cells.Item(68, 11).Value = 28000
$ws.Cells.Item(68, 12).Value = 28000
$ws.Cells.Item(68, 13).Value = 28000
$ws.Cells.Item(68, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(68, 16).Value = 1120
$ws.Cells.Item(69, 4).Value = 44676
$ws.Cells.Item(69, 11).Value = 23000
$ws.Cells.Item(69, 12).Value = 24000
$ws.Cells.Item(69, 13).Value = 23479
$ws.Cells.Item(69, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(69, 16).Value = 939
$ws.Cells.Item(70, 4).Value = 44473
$ws.Cells.Item(70, 10).Value = 85
$ws.Cells.Item(70, 11).Value = 35000
$ws.Cells.Item(70, 12).Value = 36000
$ws.Cells.Item(70, 13).Value = 35471
$ws.Cells.Item(70, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(70, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(70, 16).Value = 1419
$ws.Cells.Item(71, 4).Value = 44412
$ws.Cells.Item(71, 10).Value = 35
$ws.Cells.Item(71, 11).Value = 24000
$ws.Cells.Item(71, 12).Value = 24000
$ws.Cells.Item(71, 13).Value = 24000
$ws.Cells.Item(71, 16).Value = 960
$ws.Cells.Item(72, 4).Value = 44560
$ws.Cells.Item(72, 10).Value = 50
$ws.Cells.Item(72, 11).Value = 27000
$ws.Cells.Item(72, 12).Value = 28000
$ws.Cells.Item(72, 13).Value = 27500
$ws.Cells.Item(72, 16).Value = 1100
$ws.Cells.Item(73, 4).Value = 44578
$ws.Cells.Item(73, 11).Value = 18000
$ws.Cells.Item(73, 12).Value = 19000
$ws.Cells.Item(73, 13).Value = 18521
$ws.Cells.Item(73, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(73, 15).Value = "Provincia de Talca"
$ws.Cells.Item(73, 16).Value = 741
$ws.Cells.Item(74, 4).Value = 44571
$ws.Cells.Item(74, 10).Value = 73
$ws.Cells.Item(74, 11).Value = 15000
$ws.Cells.Item(74, 12).Value = 16000
$ws.Cells.Item(74, 13).Value = 15479
$ws.Cells.Item(74, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(74, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(74, 16).Value = 619
$ws.Cells.Item(75, 4).Value = 45085
$ws.Cells.Item(75, 10).Value = 73
$ws.Cells.Item(75, 11).Value = 30000
$ws.Cells.Item(75, 12).Value = 31000
$ws.Cells.Item(75, 13).Value = 30479
$ws.Cells.Item(75, 16).Value = 1219
$ws.Cells.Item(76, 4).Value = 44677
$ws.Cells.Item(76, 10).Value = 65
$ws.Cells.Item(76, 11).Value = 22000
$ws.Cells.Item(76, 12).Value = 23000
$ws.Cells.Item(76, 13).Value = 22462
$ws.Cells.Item(76, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(76, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(76, 16).Value = 898
$ws.Cells.Item(77, 4).Value = 44376
$ws.Cells.Item(77, 10).Value = 38
$ws.Cells.Item(77, 11).Value = 27000
$ws.Cells.Item(77, 12).Value = 27000
$ws.Cells.Item(77, 13).Value = 27000
$ws.Cells.Item(77, 16).Value = 1080
$ws.Cells.Item(78, 4).Value = 44901
$ws.Cells.Item(78, 10).Value = 65
$ws.Cells.Item(78, 11).Value = 18000
$ws.Cells.Item(78, 12).Value = 19000
$ws.Cells.Item(78, 13).Value = 18462
$ws.Cells.Item(78, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(78, 15).Value = "Región Metropolitana"
$ws.Cells.Item(78, 16).Value = 738
$ws.Cells.Item(79, 4).Value = 45111
$ws.Cells.Item(79, 10).Value = 35
$ws.Cells.Item(79, 11).Value = 24000
$ws.Cells.Item(79, 12).Value = 24000
$ws.Cells.Item(79, 13).Value = 24000
$ws.Cells.Item(79, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(79, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(79, 16).Value = 960
$ws.Cells.Item(80, 4).Value = 44511
$ws.Cells.Item(80, 10).Value = 73
$ws.Cells.Item(80, 12).Value = 17000
$ws.Cells.Item(80, 13).Value = 16479
$ws.Cells.Item(80, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(80, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(80, 16).Value = 659
$ws.Cells.Item(81, 4).Value = 44468
$ws.Cells.Item(81, 10).Value = 65
$ws.Cells.Item(81, 12).Value = 25000
$ws.Cells.Item(81, 13).Value = 24538
$ws.Cells.Item(81, 16).Value = 982
$ws.Cells.Item(82, 4).Value = 44253
$ws.Cells.Item(82, 10).Value = 38
$ws.Cells.Item(82, 11).Value = 18000
$ws.Cells.Item(82, 12).Value = 18000
$ws.Cells.Item(82, 13).Value = 18000
$ws.Cells.Item(82, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(82, 15).Value = "Provincia de Talca"
$ws.Cells.Item(82, 16).Value = 720
$ws.Cells.Item(83, 4).Value = 44634
$ws.Cells.Item(83, 11).Value = 25000
$ws.Cells.Item(83, 12).Value = 25000
$ws.Cells.Item(83, 13).Value = 25000
$ws.Cells.Item(83, 15).Value = "Provincia de Talca"
$ws.Cells.Item(83, 16).Value = 1000
$ws.Cells.Item(84, 4).Value = 44859
$ws.Cells.Item(84, 10).Value = 35
$ws.Cells.Item(84, 11).Value = 24000
$ws.Cells.Item(84, 12).Value = 24000
$ws.Cells.Item(84, 13).Value = 24000
$ws.Cells.Item(84, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(84, 16).Value = 960
$ws.Cells.Item(85, 4).Value = 44637
$ws.Cells.Item(85, 11).Value = 25000
$ws.Cells.Item(85, 12).Value = 25000
$ws.Cells.Item(85, 13).Value = 25000
$ws.Cells.Item(85, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(85, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(85, 16).Value = 1000

# Add new rows 86 and 87
$ws.Cells.Item(86, 1).Value = 3
$ws.Cells.Item(86, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(86, 3).Value = "Coquimbo"
$ws.Cells.Item(86, 4).Value = 45112
$ws.Cells.Item(86, 5).Value = 5
$ws.Cells.Item(86, 6).Value = 100112022
$ws.Cells.Item(86, 7).Value = "Arveja Verde"
$ws.Cells.Item(86, 8).Value = "Perfection"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 35
$ws.Cells.Item(86, 11).Value = 24000
$ws.Cells.Item(86, 12).Value = 24000
$ws.Cells.Item(86, 13).Value = 24000
$ws.Cells.Item(86, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(86, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(86, 16).Value = 960
$ws.Cells.Item(86, 17).Value = 25
$ws.Cells.Item(86, 18).Value = "Hortaliza"
$ws.Cells.Item(87, 1).Value = 3
$ws.Cells.Item(87, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(87, 3).Value = "Coquimbo"
$ws.Cells.Item(87, 4).Value = 44481
$ws.Cells.Item(87, 5).Value = 5
$ws.Cells.Item(87, 6).Value = 100112022
$ws.Cells.Item(87, 7).Value = "Arveja Verde"
$ws.Cells.Item(87, 8).Value = "Perfection"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 63
$ws.Cells.Item(87, 11).Value = 22000
$ws.Cells.Item(87, 12).Value = 23000
$ws.Cells.Item(87, 13).Value = 22476
$ws.Cells.Item(87, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(87, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(87, 16).Value = 899
$ws.Cells.Item(87, 17).Value = 25
$ws.Cells.Item(87, 18).Value = "Hortaliza"

# Apply date style/number format to D column of new rows to match existing date cells
$ws.Cells.Item(86, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item(87, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
